# Add columns I (I0) and J (IF) to the sheet, mirroring the existing
# header/body layout (H column = "IP") found in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1), matching style of existing headers (e.g. H1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data values for I2:J71 ---
$iVals = @(2,5,3,6,8,7,9,7,6,8,7,7,5,7,7,9,8,7,7,8,8,8,6,5,9,9,7,8,8,6,8,8,6,6,8,8,7,7,9,7,8,7,8,8,7,8,7,8,8,6,8,7,6,8,9,7,6,8,8,7,6,2,5,5,5,1,6,4,3,7)
$jVals = @(2,6,3,7,8,7,9,7,7,8,8,8,6,8,7,9,8,7,7,8,8,8,7,7,9,9,7,8,8,7,8,8,7,6,8,8,7,7,9,7,8,8,8,9,7,8,7,8,8,7,8,7,7,9,9,7,7,8,8,7,6,3,5,6,5,2,6,4,3,7)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}

$wb.Save()
